$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new events (e057, e058, e059) are inserted as new rows 87-89,
# pushing the existing e060.. rows down by three.
$ws.Rows("87:89").Insert()

# Column A: short event codes
$ws.Range("A87").Value = 'e057'
$ws.Range("A88").Value = 'e058'

# Column B: full rich-text descriptions
$ws.Range("B87").Value = '<Bold>e057 Fire 2 Inch Mortar</Bold> 
<InlineUIContainer><Button Content=''r4.74.4'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>   
<InlineUIContainer><Button Content=''r18.12'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
One Smoke Marker is placed in the close zone in front of turret, and mark off one smoke grenade on the After Action Report 
<InlineUIContainer><Button Content=''AAR'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name=''c58LFireMortar'' Height=''100'' Width=''100''></Image></InlineUIContainer>  '
$ws.Range("B88").Value = '<Bold>e058 Throw Grenade</Bold> 
<InlineUIContainer><Button Content=''r4.74.4'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>   
<InlineUIContainer><Button Content=''r18.11'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
One Smoke Marker is placed your tank by a crewman throwing a grenade. Mark off one smoke grenade on the After Action Report 
<InlineUIContainer><Button Content=''AAR'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name=''c70ThrowSmokeGrenade'' Height=''100'' Width=''100''></Image></InlineUIContainer>  '

$ws.Range("A89").Value = 'e059'
$ws.Range("B89").Value = '<Bold>e059 Restock Ready Rack</Bold> 
<InlineUIContainer><Button Content=''r4.74.4'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>   
<InlineUIContainer><Button Content=''r16.23'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Loader restocks the ready rack from available ammo shown on the After Action Report 
<InlineUIContainer><Button Content=''AAR'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>. Click buttons to adjust or the image when done.
<LineBreak/><LineBreak/>
   <InlineUIContainer><Button Name=''HeMinus'' Content=''   -   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>  
<InlineUIContainer><Button Name=''HePlus''  Content=''   +   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> for HE Ammo<LineBreak/>
   <InlineUIContainer><Button Name=''ApMinus'' Content=''   -   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>  
<InlineUIContainer><Button Name=''ApPlus'' Content=''   +   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> for AP Ammo<LineBreak/>
   <InlineUIContainer><Button Name=''WpMinus'' Content=''   -   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>  
<InlineUIContainer><Button Name=''WpPlus'' Content=''   +   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> for WP Ammo<LineBreak/>
   <InlineUIContainer><Button Name=''HcbiMinus'' Content=''   -   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>  
<InlineUIContainer><Button Name=''HcbiPlus'' Content=''   +   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> for HCBI Ammo<LineBreak/>
   <InlineUIContainer><Button Name=''HvapMinus'' Content=''   -   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>  
<InlineUIContainer><Button Name=''HvapPlus'' Content=''   +   '' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> for HVAP Ammo<LineBreak/>
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name=''c60LRestockReadyRack'' Height=''100'' Width=''100''></Image></InlineUIContainer>  '

# Row heights to fit the new wrapped text content
$ws.Rows(87).RowHeight = 120
$ws.Rows(88).RowHeight = 120
$ws.Rows(89).RowHeight = 285

# Scroll/select the newly added area, as in the saved view
$excel.ActiveWindow.ScrollRow = 89
$ws.Range("B89").Select()

